$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-15 (A: date serial, B: TISG, C: fcs, E: MYDIR)
# Columns D (buy) and F (need_to_buy) are cleared.
$data = @(
    @(46049, 11802.7698759612, 11653.0209379243, 7638.1684028183),
    @(46050, 12147.9154592183, 11481.903261406,  8146.83256910503),
    @(46051, 12043.293228742,  11580.3752723273, 8074.34442205525),
    @(46052, 12126.1048736501, 11087.2291643257, 8150.02269995864),
    @(46053, 4757.35041428795, 7833.31496487201, 7825.52543964492),
    @(46054, 5050.65013921885, 7896.54178192751, 7899.72962513176),
    @(46055, 11455.6405730248, 11045.1361485746, 7815.50065662757),
    @(46056, 11455.6405730248, 10853.9880842006, 7815.50065662757),
    @(46057, 11455.6405730248, 10799.4114611398, 7815.50065662757),
    @(46058, 11455.6405730248, 10792.5346956402, 7815.50065662757),
    @(46059, 11455.6405730248, 10182.0211295624, 7815.50065662757),
    @(46060, 4820.25116453669, 7341.76983611613, 7442.14865737584),
    @(46061, 4720.83339186724, 7345.94629602197, 7433.4753349264),
    @(46062, 11066.9810198305, 10627.2921415579, 7339.15639816247)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).ClearContents()
    $ws.Cells.Item($row, 5).Value = $entry[3]
    $ws.Cells.Item($row, 6).ClearContents()
    $row++
}
